$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    $rng.Text = $newText
}

Replace-Text "Wanahisabati wanaocheza:" "The playful mathematicians:"
Replace-Text "** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino" "** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino"
Replace-Text "[Muziki]" "[Music]"
Replace-Text "kuna wanahisabati wawili, tupige simu" "there are two mathematicians, let's call"
Replace-Text "Fil na Mike wanaokutana" "them Fil and Mike who meet each other"
Replace-Text "tena baada ya muda mrefu. Baada ya baadhi" "again after a long time. After some"
Replace-Text "kuzungumza, Phil anasema ana watoto watatu, basi" "chatting, Phil says he has three children, then"
Replace-Text "Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil," "Mike, astonished, asks: 'How old are they?' Fil,"
Replace-Text "kuwa mwanahisabati mchezaji, anajibu" "being a playful mathematician, answers"
Replace-Text "'Wewe niambie! Nitakupa kidokezo: ikiwa wewe" "'You tell me! I'll give you a hint: if you"
Replace-Text "zidisheni enzi tatu pamoja ninyi" "multiply the three ages together you"
Replace-Text "pata 36.' Mike huchukua wakati mwingine kufikiria" "get 36.' Mike takes sometimes to think"
Replace-Text "na kusema: 'Samahani Fil, lakini nahitaji" "and says: 'I'm sorry Fil, but I do need"
Replace-Text "kidokezo kingine. Kwa hivyo Fil anamwambia Mike:" "another hint. So Fil tells Mike:"
Replace-Text "'Ndiyo, hakika, hapa ni: kama alikuwa na hadi" "'Yes, sure, here it is: if you had up to"
Replace-Text "miaka mitatu unapata idadi ya hesabu" "three ages you get the number of math"
Replace-Text "karatasi tunachapisha pamoja. Je, unaikumbuka?'" "papers we publish together. Do you remember it?'"
Replace-Text "'Ndio nakumbuka wangapi, lakini bado" "'Yes I do remember How many, but still"
Replace-Text "Sina taarifa za kutosha! nahitaji" "I do not have enough information! I need"
Replace-Text "angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo" "at least one more.' Fil says: 'Yes don't"
Replace-Text "wasiwasi lakini hii ni ya mwisho:" "worry but this is the last one:"
Replace-Text "Mdogo ana macho ya blues.' Na" "The youngest one has blues eyes.' And"
Replace-Text "ghafla Mike anapata jibu. Wewe" "suddenly Mike gets the answer. You"
Replace-Text "sikia mazungumzo lakini hujui" "hear the conversation but you don't know"
Replace-Text "ni karatasi ngapi walichapisha pamoja." "how many papers they published together."
Replace-Text "Hata hivyo, unataka kujua umri wa" "However, you do want to know the ages of"
Replace-Text "watoto watatu. Je, unaweza kuwahesabu" "the three children. Can you figure them"
Replace-Text "nje?" "out?"
Replace-Text "[Muziki]" "[Music]"
